$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (shared, rich-text strings) ---
# A8: "Volume 31   Number  34" -> "...35"  (replace the trailing "34" run)
$hdrVol = $ws.Range("A8")
$ch = $hdrVol.Characters(21, 2)
$ch.Text = "35"

# C9: "Report Covering the Week  8/19/2024  Through  8/25/2024"
#     -> "...8/26/2024  Through  9/1/2024"
$hdrWeek = $ws.Range("C9")
$chStart = $hdrWeek.Characters(27, 9)
$chStart.Text = "8/26/2024"
$chEnd = $hdrWeek.Characters(47, 9)
$chEnd.Text = "9/1/2024"

# --- Crime-complaints table updates (rows 14-33) ---
$ws.Range("N14").Value = -94.117647058823
$ws.Range("C15").Value = 1
$ws.Range("F17").Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 1
$ws.Range("F17").Copy()
$ws.Range("F15").PasteSpecial(-4122)
$ws.Range("G15").Value = 3
$ws.Range("H15").Value = -66.666666666666
$ws.Range("I15").Value = 17
$ws.Range("J15").Value = 15
$ws.Range("K15").Value = 13.333333333333
$ws.Range("L15").Value = 30.769230769230
$ws.Range("M15").Value = 21.428571428571
$ws.Range("N15").Value = -5.555555555555
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 4
$ws.Range("E16").Value = 25
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 17
$ws.Range("H16").Value = 0
$ws.Range("I16").Value = 142
$ws.Range("J16").Value = 146
$ws.Range("K16").Value = -2.739726027397
$ws.Range("L16").Value = 6.766917293233
$ws.Range("M16").Value = -21.546961325966
$ws.Range("N16").Value = -77.207062600321
$ws.Range("C17").Value = 6
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 200
$ws.Range("G17").Value = 8
$ws.Range("H17").Value = 162.5
$ws.Range("I17").Value = 203
$ws.Range("J17").Value = 183
$ws.Range("K17").Value = 10.928961748633
$ws.Range("L17").Value = 2.525252525252
$ws.Range("M17").Value = 32.679738562091
$ws.Range("N17").Value = -1.932367149758
$ws.Range("D18").Value = 6
$ws.Range("E18").Value = -66.666666666666
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 19
$ws.Range("H18").Value = -21.052631578947
$ws.Range("I18").Value = 152
$ws.Range("J18").Value = 129
$ws.Range("K18").Value = 17.829457364341
$ws.Range("L18").Value = -17.837837837837
$ws.Range("M18").Value = -50.326797385620
$ws.Range("N18").Value = -89.064748201438
$ws.Range("C19").Value = 9
$ws.Range("D19").Value = 18
$ws.Range("E19").Value = -50
$ws.Range("F19").Value = 42
$ws.Range("G19").Value = 63
$ws.Range("H19").Value = -33.333333333333
$ws.Range("I19").Value = 445
$ws.Range("J19").Value = 463
$ws.Range("K19").Value = -3.887688984881
$ws.Range("L19").Value = 2.064220183486
$ws.Range("M19").Value = 60.649819494584
$ws.Range("N19").Value = 0.678733031674
$ws.Range("C20").Value = 13
$ws.Range("D20").Value = 11
$ws.Range("E20").Value = 18.181818181818
$ws.Range("F20").Value = 40
$ws.Range("G20").Value = 43
$ws.Range("H20").Value = -6.976744186046
$ws.Range("I20").Value = 286
$ws.Range("J20").Value = 247
$ws.Range("K20").Value = 15.789473684210
$ws.Range("L20").Value = 52.941176470588
$ws.Range("M20").Value = 10
$ws.Range("N20").Value = -88.240131578947
$ws.Range("C21").Value = 36
$ws.Range("D21").Value = 42
$ws.Range("E21").Value = -14.285714285714
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 153
$ws.Range("H21").Value = -11.111111111111
$ws.Range("I21").Value = 1246
$ws.Range("J21").Value = 1186
$ws.Range("K21").Value = 5.059021922428
$ws.Range("L21").Value = 7.972270363951
$ws.Range("M21").Value = 4.530201342281
$ws.Range("N21").Value = -75.706765451355
$ws.Range("G22").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("G22").PasteSpecial(-4122)
$ws.Range("H22").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("H22").PasteSpecial(-4122)
$ws.Range("M22").Value = 8.333333333333
$ws.Range("C24").Value = 38
$ws.Range("D24").Value = 44
$ws.Range("E24").Value = -13.636363636363
$ws.Range("F24").Value = 114
$ws.Range("G24").Value = 132
$ws.Range("H24").Value = -13.636363636363
$ws.Range("I24").Value = 976
$ws.Range("J24").Value = 901
$ws.Range("K24").Value = 8.324084350721
$ws.Range("L24").Value = -0.102354145342
$ws.Range("M24").Value = 34.620689655172
$ws.Range("C25").Value = 22
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = 22.222222222222
$ws.Range("F25").Value = 61
$ws.Range("G25").Value = 52
$ws.Range("H25").Value = 17.307692307692
$ws.Range("I25").Value = 384
$ws.Range("J25").Value = 315
$ws.Range("K25").Value = 21.904761904761
$ws.Range("L25").Value = -3.517587939698
$ws.Range("C26").Value = 14
$ws.Range("D26").Value = 12
$ws.Range("E26").Value = 16.666666666666
$ws.Range("F26").Value = 46
$ws.Range("G26").Value = 34
$ws.Range("H26").Value = 35.294117647058
$ws.Range("I26").Value = 405
$ws.Range("J26").Value = 322
$ws.Range("K26").Value = 25.776397515528
$ws.Range("L26").Value = 14.730878186968
$ws.Range("M26").Value = -18.511066398390
$ws.Range("C27").Value = 1
$ws.Range("F17").Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("E27").Value = 0
$ws.Range("F27").Value = 1
$ws.Range("F17").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = -66.666666666666
$ws.Range("I27").Value = 23
$ws.Range("J27").Value = 24
$ws.Range("K27").Value = -4.166666666666
$ws.Range("L27").Value = 4.545454545454
$ws.Range("C28").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("D28").Value = "'0"
$ws.Range("C23").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("E28").Value = "***.*"
$ws.Range("C23").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("F28").Value = 2
$ws.Range("G28").Value = 4
$ws.Range("L28").Value = -57.692307692307
$ws.Range("N29").Value = -91.304347826087
$ws.Range("N30").Value = -90
$ws.Range("C33").Value = 1
$ws.Range("F17").Copy()
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("F33").Value = 1
$ws.Range("F17").Copy()
$ws.Range("F33").PasteSpecial(-4122)
$ws.Range("G33").Value = 1
$ws.Range("H33").Value = 0
$ws.Range("I33").Value = 2
$ws.Range("K33").Value = -66.666666666666
$ws.Range("L33").Value = -50

